$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 159
$ws.Range("F3").Value = 204
$ws.Range("F6").Value = 1351
$ws.Range("F7").Value = 75
$ws.Range("F9").Value = 397
$ws.Range("F10").Value = 456
$ws.Range("F11").Value = 828
$ws.Range("F12").Value = 230
$ws.Range("F13").Value = 752
$ws.Range("F14").Value = 324
$ws.Range("F15").Value = 490
$ws.Range("F16").Value = 95
$ws.Range("F17").Value = 1061
$ws.Range("F18").Value = 510
$ws.Range("F19").Value = 294
$ws.Range("F20").Value = 424
$ws.Range("F21").Value = 112
$ws.Range("F22").Value = 238
$ws.Range("F23").Value = 32
$ws.Range("F24").Value = 58
$ws.Range("F25").Value = 497
$ws.Range("F26").Value = 463
$ws.Range("F28").Value = 332

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 54
$ws.Range("F10").Value = 634
$ws.Range("F12").Value = 156

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 159
$ws.Range("F5").Value = 204
$ws.Range("F8").Value = 1351
$ws.Range("F10").Value = 75
$ws.Range("F13").Value = 54
$ws.Range("F14").Value = 397
$ws.Range("F17").Value = 456
$ws.Range("F18").Value = 828
$ws.Range("F19").Value = 230
$ws.Range("F20").Value = 752
$ws.Range("F21").Value = 324
$ws.Range("F22").Value = 490
$ws.Range("F23").Value = 95
$ws.Range("F24").Value = 1061
$ws.Range("F25").Value = 510
$ws.Range("F28").Value = 294
$ws.Range("F29").Value = 424
$ws.Range("F30").Value = 634
$ws.Range("F31").Value = 112
$ws.Range("F33").Value = 238
$ws.Range("F34").Value = 32
$ws.Range("F35").Value = 58
$ws.Range("F36").Value = 156
$ws.Range("F38").Value = 497
$ws.Range("F41").Value = 463
$ws.Range("F43").Value = 332
